# Update countries & provincias Spain
#
# This script refreshes the COVID-19 stats table on the "Pais" sheet.
# Three country names change position relative to their neighbours in
# the underlying data (the sheet is sorted by total cases, so a country
# overtaking/falling behind another effectively swaps which row holds
# which name) and a handful of rows receive refreshed numeric figures
# (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    # Rows 28-30: Bielorrusia overtakes Singapur/Irlanda and gets fresh
    # data; Singapur/Irlanda inherit the figures that used to belong to
    # the row above them.
    @{ Row=28;  Country="Bielorrusia";    B=23906; C=933; D=6531;  E=17240; F=92; G=4; H=135 },
    @{ Row=29;  Country="Singapur";       B=23822; C=486; D=2721;  E=21081; F=23; G=0; H=20 },
    @{ Row=30;  Country="Irlanda";        B=22996; C=0;   D=17110; E=4428;  F=72; G=0; H=1458 },

    # Row 53: Australia - figures refreshed, name unchanged.
    @{ Row=53;  Country="Australia";      B=6948;  C=7;   D=6179;  E=672;   F=16; G=0; H=97 },

    # Row 59: Kazajistan - figures refreshed, name unchanged.
    @{ Row=59;  Country="Kazajistan";     B=5138;  C=48;  D=1941;  E=3166;  F=31; G=0; H=31 },

    # Row 89: Eslovenia - figures refreshed, name unchanged.
    @{ Row=89;  Country="Eslovenia";      B=1460;  C=3;   D=256;   E=1102;  F=10; G=0; H=102 },

    # Rows 140-141: Etiopia overtakes Cabo Verde and gets fresh data;
    # Cabo Verde inherits the figures that used to belong to the row
    # above it.
    @{ Row=140; Country="Etiopia";        B=250;   C=11;  D=105;   E=140;   F=1;  G=0; H=5 },
    @{ Row=141; Country="Cabo Verde";     B=246;   C=0;   D=56;    E=188;   F=0;  G=0; H=2 },

    # Rows 192-193: Belice overtakes Nueva Caledonia; only Casos activos
    # (D) and Muertes (H) differ between the two, the rest coincide.
    @{ Row=192; Country="Belice";          B=18;   C=0;   D=16;    E=0;     F=0;  G=0; H=2 },
    @{ Row=193; Country="Nueva Caledonia"; B=18;   C=0;   D=18;    E=0;     F=0;  G=0; H=0 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Country
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}
